$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.455362044514542
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 3.537761648806719
$ws.Range("E2").Value = 10.19245300693656
$ws.Range("G2").Value = 16.84135478251809

$ws.Range("B3").Value = 0.00001292064567892659
$ws.Range("C3").Value = 0.00006240767534437808
$ws.Range("D3").Value = 22.3905356188092
$ws.Range("E3").Value = 1133.036916526867
$ws.Range("G3").Value = 1155.427527473998
